# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45177 (2023-09-08) to 45178 (2023-09-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 185; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}
